$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.512.24'
$ws.Range('E2').Value = '  +1.62%  '

$ws.Range('D3').Value = '3.742.16'
$ws.Range('E3').Value = '  -0.31%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.83'
$ws.Range('E5').Value = '  -0.13%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.46'
$ws.Range('E6').Value = '  -1.12%  '

$ws.Range('D7').Value = '3.740.59'
$ws.Range('E7').Value = '  -0.34%  '

$ws.Range('E8').Value = '  -0.06%  '

$ws.Range('E9').Value = '  -0.80%  '

$ws.Range('E10').Value = '  -3.20%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.49'
$ws.Range('E11').Value = '  +0.17%  '

$ws.Range('E12').Value = '  -0.71%  '

$ws.Range('E13').Value = '  -6.02%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.33'
$ws.Range('E14').Value = '  -0.38%  '

$ws.Range('D15').Value = '4.374.32'
$ws.Range('E15').Value = '  -0.18%  '

$ws.Range('D16').Value = '3.743.11'
$ws.Range('E16').Value = '  -0.28%  '

$ws.Range('D17').Value = '68.477.04'
$ws.Range('E17').Value = '  +1.68%  '

$ws.Range('E18').Value = '  -3.44%  '

$ws.Range('E19').Value = '  -2.16%  '

$ws.Range('E20').Value = '  -0.25%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.69'
$ws.Range('E21').Value = '  +1.67%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '467.80'
$ws.Range('E22').Value = '  +0.30%  '

$ws.Range('E23').Value = '  -2.50%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.24'
$ws.Range('E24').Value = '  +0.75%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000144'
$ws.Range('E25').Value = '  -2.19%  '

$ws.Range('E26').Value = '  -0.18%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.08'
$ws.Range('E27').Value = '  -0.53%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.15'
$ws.Range('E28').Value = '  -1.41%  '

$ws.Range('E29').Value = '  -0.10%  '

$ws.Range('D30').Value = '3.891.99'
$ws.Range('E30').Value = '  -0.28%  '

$ws.Range('E31').Value = '  -3.88%  '

$ws.Range('E32').Value = '  -3.99%  '

$ws.Range('E33').Value = '  -1.65%  '

$ws.Range('E34').Value = '  -1.74%  '

$ws.Range('E35').Value = '  +1.68%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'

$ws.Range('D37').Value = '3.701.17'
$ws.Range('E37').Value = '  -0.37%  '

$ws.Range('E38').Value = '  -1.45%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.46'
$ws.Range('E39').Value = '  -8.73%  '

$ws.Range('E40').Value = '  +0.89%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').Value = '  +0.02%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.80'
$ws.Range('E42').Value = '  -0.24%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.07%  '

$ws.Range('E45').Value = '  -1.85%  '

$ws.Range('E46').Value = '  -0.20%  '

$ws.Range('E47').Value = '  -0.99%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '43.05'
$ws.Range('E48').Value = '  +10.44%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '45.90'
$ws.Range('E49').Value = '  +0.15%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '146.99'
$ws.Range('E50').Value = '  +5.76%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '391.35'
$ws.Range('E51').Value = '  -1.59%  '
